$d = $word.ActiveDocument

# 1) "Servicio al cliente a bordo" -> "Personal a bordo"
$d.Content.Find.Execute("Servicio al cliente a bordo", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Personal a bordo", 2)

# Locate the paragraph that now holds the edited text so the bookmark
# move below doesn't depend on a hard-coded paragraph index.
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Range.Text -like "Personal a bordo*") {
        $target = $candidate
    }
}

$pEnd = $target.Range.End

# Move the document's "_GoBack" bookmark (originally sitting after
# "Control de usuarios") so it becomes an empty bookmark right after the
# new "Personal a bordo" run. Adding a bookmark with a name that already
# exists replaces/relocates the existing one, so this single Add() both
# removes the old occurrence and creates the new one.
#
# A zero-length range placed exactly on the paragraph-mark position would
# land in the wrong spot, so temporarily insert a marker character right
# before the paragraph mark, anchor the bookmark next to the real text,
# then delete the marker again -- leaving the bookmark collapsed
# immediately after "bordo".
$tail = $d.Range($pEnd - 1, $pEnd - 1)
$tail.InsertAfter("X")

$bmPos = $d.Range($pEnd - 1, $pEnd - 1)
$d.Bookmarks.Add("_GoBack", $bmPos)

$marker = $d.Range($pEnd - 1, $pEnd)
$marker.Delete()

# 2) Collapse the three runs "Lugares "/"Salida"/"/Destino" into a single
#    run "Lugares Salida/Destino".
$d.Content.Find.Execute("Lugares Salida/Destino", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Lugares Salida/Destino", 2)
